# Update the "2018 World Cup" results for the matches played through June 22.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("2018 World Cup")

# Group stage match results (matchday played on/around Jun 21-22, 2018)
$ws.Range("F27").Value = 1   # Francia
$ws.Range("G27").Value = 0   # Peru

$ws.Range("F28").Value = 1   # Dinamarca
$ws.Range("G28").Value = 1   # Australia

$ws.Range("F29").Value = 0   # Argentina
$ws.Range("G29").Value = 3   # Croacia

$ws.Range("F30").Value = 2   # Nigeria
$ws.Range("G30").Value = 0   # Islandia

$ws.Range("F31").Value = 2   # Brasil
$ws.Range("G31").Value = 0   # Costa Rica

$ws.Range("F32").Value = 1   # Serbia
$ws.Range("G32").Value = 2   # Suiza

# Restore the view / selection state as it was left after the edit
$ws.Application.Goto($ws.Range("B25"), $false)
$ws.Range("J35").Select()
